$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 3-8 (columns D, I, J, K, L, M, N, O, P, Q).
# This represents the rows being re-sorted (e.g. by date), with each
# row taking on the values previously held by another row.
$rows = @{
    3 = @{ D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; O = "Región de Arica y Parinacota"; P = 1500; Q = 20 }
    4 = @{ D = 44315; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel";    O = "Región de Arica y Parinacota"; P = 1000; Q = 15 }
    5 = @{ D = 44280; I = "Primera";  J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; O = "Región de Arica y Parinacota"; P = 1389; Q = 18 }
    6 = @{ D = 44285; I = "Primera";  J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; O = "Región de Arica y Parinacota"; P = 1389; Q = 18 }
    7 = @{ D = 44313; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; O = "Región de Arica y Parinacota"; P = 1000; Q = 15 }
    8 = @{ D = 44313; I = "Primera";  J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; O = "Región de Arica y Parinacota"; P = 1500; Q = 20 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
